$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-60 down to 41-61
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44977
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100101
$ws.Cells.Item(40, 8).Value = "Berries"
$ws.Cells.Item(40, 9).Value = 100101001
$ws.Cells.Item(40, 10).Value = "Arándano (blue)"
$ws.Cells.Item(40, 11).Value = "Sin especificar"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 200
$ws.Cells.Item(40, 14).Value = 2000
$ws.Cells.Item(40, 15).Value = 2200
$ws.Cells.Item(40, 16).Value = 2100
$ws.Cells.Item(40, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(40, 19).Value = 1050
$ws.Cells.Item(40, 20).Value = 2
